$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7971
$ws1.Range("F5").Value = 5824
$ws1.Range("F7").Value = 84
$ws1.Range("F10").Value = 282
$ws1.Range("F11").Value = 352

# Sheet "全部类型" (4th sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7971
$ws4.Range("F5").Value = 5824
$ws4.Range("F7").Value = 84
$ws4.Range("F10").Value = 282
$ws4.Range("F14").Value = 352
